# Add a "PO Forecast" sheet populated with forecast data, and rename the
# "Requested quantity" header on the existing "Weekly Quantity" / "Monthly
# Trend" sheets to their new, more specific names.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers -------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet as the last (3rd) sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row --------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) -------------
$dataBlob = @"
44934.99999999999|345|-19.5510030073467|688.2218698717106
44941.99999999999|342|8.510782297323868|691.0980072117333
44948.99999999999|340|-20.31966872611999|686.3270277381719
44955.99999999999|338|-7.323274820908957|685.4753203139096
44962.99999999999|335|12.63211517214473|685.7699672253829
44969.99999999999|333|-17.99165880893628|667.1170488252986
44976.99999999999|331|-22.38887021641529|669.7286368096466
44983.99999999999|328|-26.49024282724122|679.4272011137825
44990.99999999999|326|-36.84017872875516|677.9462782283488
44997.99999999999|324|-24.73953175888536|663.0219411876814
45018.99999999999|317|-10.83848209935677|647.9218004148671
45060.99999999999|303|-35.16384496026387|638.7754757032565
45067.99999999999|300|-54.22643367617699|620.609478289151
45074.99999999999|298|-50.7824937248345|658.2303356573709
45081.99999999999|296|-47.38301485677529|636.099494348969
45088.99999999999|293|-75.12203912483047|652.1470163712273
45095.99999999999|291|-75.89184045778133|626.3447430941609
45102.99999999999|289|-31.76370624190164|665.2322979662561
45109.99999999999|286|-56.99669387866615|624.3351397925982
45116.99999999999|284|-58.47600451779857|656.7832882578225
45123.99999999999|282|-72.56528532828105|621.3864486043724
45137.99999999999|277|-58.88852438743053|672.9078333214946
45144.99999999999|275|-72.99242096607676|609.6489405930225
45158.99999999999|270|-91.20582919071835|614.7855425227057
45165.99999999999|268|-83.51019347134407|613.0331092777222
45172.99999999999|265|-77.51326738306278|591.2315946943029
45179.99999999999|263|-85.68749804698612|601.6338371079014
45193.99999999999|258|-93.39418558288745|596.2122583380341
45200.99999999999|256|-85.0699769361165|630.0437527775626
45207.99999999999|254|-101.6202177114198|595.9513445563022
45228.99999999999|247|-121.2018536369302|565.4962411625565
45242.99999999999|242|-108.4157484787357|605.5120494206094
45249.99999999999|240|-120.0695449997605|565.6838344673314
45256.99999999999|238|-98.1569368975275|596.3408382910482
45270.99999999999|233|-123.1507567079288|558.8030000593279
45277.99999999999|231|-99.07937779528889|583.3245498068366
45298.99999999999|224|-127.3021344863284|572.5192337078789
45305.99999999999|221|-101.7460589024874|574.0370140772895
45312.99999999999|219|-131.1424980051536|550.552178025153
45340.99999999999|210|-147.0474406954485|547.3779971616149
45347.99999999999|207|-151.7080674280936|539.7203412191986
45354.99999999999|205|-126.007092009298|555.0486016613869
45368.99999999999|200|-166.3831988910805|532.5812654031087
45375.99999999999|198|-130.7493440461967|534.8309792964264
45417.99999999999|184|-160.9103694855962|504.680275581601
45438.99999999999|177|-156.5322544396114|504.4191748126494
45473.99999999999|165|-191.3342080315921|500.2519034245894
45480.99999999999|163|-177.263427018036|506.7908271923731
45487.99999999999|161|-172.2422641480331|492.3691708090032
45515.99999999999|151|-224.6512425310385|481.588715308413
45522.99999999999|149|-206.3073305051263|484.0604105502857
45529.99999999999|147|-208.713621373703|501.0721262433239
45543.99999999999|142|-214.2704146040026|496.6898215604519
45550.99999999999|140|-212.3792872544751|482.3863591813515
45557.99999999999|137|-215.36396393634|482.3499601960738
45564.99999999999|135|-233.0152440822609|496.380706145786
45571.99999999999|133|-192.1506443929858|492.486980214681
45578.99999999999|130|-188.4755501528893|481.9131718326693
45585.99999999999|128|-213.3871742481033|481.3796036410936
45592.99999999999|126|-214.6214659910042|488.9176511051792
45599.99999999999|123|-202.4355411006563|483.1175582597609
"@

$dataRows = $dataBlob.Trim().Split("`n")
$rowIndex = 2
foreach ($line in $dataRows) {
    $fields = $line.Trim().Split("|")
    $wsForecast.Cells.Item($rowIndex, 1).Value = [double]$fields[0]
    $wsForecast.Cells.Item($rowIndex, 2).Value = [double]$fields[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = [double]$fields[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = [double]$fields[3]
    $rowIndex = $rowIndex + 1
}

$lastRow = $rowIndex - 1

# --- 5. Formatting: reuse the header / date styles from the Weekly sheet
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range(("A2:A" + $lastRow)).PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Host ("PO Forecast sheet populated with " + ($rowIndex - 2) + " rows")
